# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 22:14"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8269122
$ws.Range("C4").Value = 52807
$ws.Range("D4").Value = 5344380
$ws.Range("E4").Value = 2701328
$ws.Range("G4").Value = 698
$ws.Range("H4").Value = 223414

# Row 13
$ws.Range("E13").Value = 697385
$ws.Range("G13").Value = 178
$ws.Range("H13").Value = 33303

# Row 21
$ws.Range("B21").Value = 356668
$ws.Range("C21").Value = 7852
$ws.Range("D21").Value = 287600
$ws.Range("E21").Value = 59232
$ws.Range("G21").Value = 26
$ws.Range("H21").Value = 9836

# Row 50
$ws.Range("B50").Value = 94348
$ws.Range("C50").Value = 1196
$ws.Range("D50").Value = 58269
$ws.Range("E50").Value = 34911
$ws.Range("G50").Value = 9
$ws.Range("H50").Value = 1168

# Row 103
$ws.Range("B103").Value = 12215
$ws.Range("C103").Value = 112
$ws.Range("D103").Value = 10360
$ws.Range("E103").Value = 1724
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 131

# Row 118
$ws.Range("B118").Value = 7526
$ws.Range("C118").Value = 82
$ws.Range("D118").Value = 6425
$ws.Range("E118").Value = 1019
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 82

# Row 123
$ws.Range("B123").Value = 5842
$ws.Range("C123").Value = 6
$ws.Range("D123").Value = 4735
$ws.Range("E123").Value = 926

# Row 124
$ws.Range("B124").Value = 5746
$ws.Range("C124").Value = 13
$ws.Range("D124").Value = 5392
$ws.Range("E124").Value = 239

# Row 125
$ws.Range("B125").Value = 5449
$ws.Range("C125").Value = 6
$ws.Range("D125").Value = 5372
$ws.Range("E125").Value = 16

# Row 138
$ws.Range("B138").Value = 4289
$ws.Range("C138").Value = 4
$ws.Range("D138").Value = 3947
$ws.Range("E138").Value = 310
